# Calculate Arbitrage Bid / Ask
# Duplicates the existing BTC/USDT - ETH/BTC - ETH/USDT arbitrage calculation
# into a "Kaufen/Bid" block (the original calc, shifted down) and a new
# "Verkaufen/Ask" block (the mirrored calculation), with header labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing table (rows 1 & 2, originally) down by two rows so the
# original header/data/formula rows land on 3/4/6 and free up row 1 for a
# new "Base Coin / USDT" caption. References inside the existing formulas
# are updated automatically by the row insert.
$ws.Rows("1:2").Insert()

# --- New caption above the (now shifted) first table ----------------------
$ws.Range("A1").Value = "Base Coin"
$ws.Range("C1").Value = "USDT"

# --- Second block: caption -------------------------------------------------
$ws.Range("A11").Value = "Base Coin"
$ws.Range("C11").Value = "BTC"

# --- "Verkaufen" / "Ask" labels below the second block ---------------------
$ws.Range("A18").Value = "Verkaufen"
$ws.Range("A19").Value = "Ask"

# --- "Kaufen" / "Bid" labels below the first table --------------------------
$ws.Range("A8").Value = "Kaufen"
$ws.Range("A9").Value = "Bid"

# --- Second block: headers + source data (duplicate of the first table) ----
$ws.Range("A13").Value = "BTC/USDT"
$ws.Range("C13").Value = "ETH/BTC"
$ws.Range("E13").Value = "ETH/USDT"

$ws.Range("A14").Value = 54899.31
$ws.Range("C14").Value = 0.030953999999999999
$ws.Range("E14").Value = 1699.5

# --- Second block: mirrored arbitrage formulas ------------------------------
$ws.Range("A16").Formula = "=A14/1"
$ws.Range("C16").Formula = "=A16/C14"
$ws.Range("E16").Formula = "=C16*E14"

# --- Number formatting: 7-decimal thousands format on the result rows ------
$ws.Range("A6").NumberFormat = "#,##0.0000000"
$ws.Range("C6").NumberFormat = "#,##0.0000000"
$ws.Range("E6").NumberFormat = "#,##0.0000000"

$ws.Range("A16").NumberFormat = "#,##0.0000000"
$ws.Range("C16").NumberFormat = "#,##0.0000000"
$ws.Range("E16").NumberFormat = "#,##0.0000000"

# --- Column widths (best-fit approximations) --------------------------------
$ws.Columns("A").ColumnWidth = 13.26
$ws.Columns("C").ColumnWidth = 15.94
$ws.Columns("E").ColumnWidth = 19.62

# --- Page setup ---------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection as left by the author -------------------------------------------
$ws.Range("A10").Select()
